$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "N/A" placeholders in E13 and E17 with the actual computed
# portfolio-return figure. This lets the dependent H/I formulas
# (C*(E-F) and (C-D)*(E-F)) resolve to numbers instead of #VALUE! errors,
# and the SUM/SUMPRODUCT rollups in row 10 and row 14 recalc accordingly.
$ws.Range("E13").Value = -0.0194690265486725
$ws.Range("E17").Value = -0.0194690265486725

# Move the active selection, matching the saved view state.
$ws.Range("H22").Select()

$wb.Save()
